# Apply cryptos list update (price/volume refresh + two rank swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 19/20 swap rank order: Chainlink <-> Polkadot ---
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.31"
$ws.Range("E19").Value = "  -3.46%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  -2.90%  "

# --- Rows 25/26 swap rank order: WrappedeETH <-> Dai ---
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.662.81"
$ws.Range("E26").Value = "  -4.07%  "

# --- Remaining Price (D) / Volume(1h) (E) refreshed values ---
$ws.Range("D2").Value = "66.901.98"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").Value = "3.521.64"
$ws.Range("E3").Value = "  -4.05%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.28"
$ws.Range("E5").Value = "  -5.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.75"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("D7").Value = "3.518.22"
$ws.Range("E7").Value = "  -4.06%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.424"
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("E13").Value = "  -4.60%  "
$ws.Range("D14").Value = "4.114.93"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.45"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").Value = "3.522.76"
$ws.Range("E16").Value = "  -4.53%  "
$ws.Range("D17").Value = "66.912.95"
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.12"
$ws.Range("E21").Value = "  -5.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.16"
$ws.Range("E22").Value = "  -8.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.623"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.34"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").Value = "  -5.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.14"
$ws.Range("E29").Value = "  -9.15%  "
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.65"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.52"
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("D35").Value = "3.514.40"
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.09"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("E37").Value = "  -7.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.92"
$ws.Range("E38").Value = "  -5.56%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.76"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.12"
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.49"
$ws.Range("E43").Value = "  -5.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0855"
$ws.Range("E44").Value = "  -3.78%  "
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.25"
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.36"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.55"
$ws.Range("E48").Value = "  -5.03%  "
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.50"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  -4.47%  "
